$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All D-column "Price" values are plain text in the source data (thousand-separator
# dotted numbers, fixed decimal places, scientific-looking small decimals, etc.).
# Force text storage via NumberFormat "@" so Excel does not silently convert the
# assigned string into a numeric value (which would drop trailing/representative zeros),
# then restore the default "Normal" style so no stray formatting is left on the cell.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.985.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.17%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9961"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6303"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9981"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.23%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07520"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2935"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.82%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.37%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07712"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.98%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.834.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.30%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.984"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6699"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.80%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.78%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009644"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.16%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.074"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.023.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.76%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.56%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "226.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.84%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9974"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.33%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.158"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9982"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1403"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.50%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.514"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "

# Row 28
$ws.Range("E28").Value = "  -0.12%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.119"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.77%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.062"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.195"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.29%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05366"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.27%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.857"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.28%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7429"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.47%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.137"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.36%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.648"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.43%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.243.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.67%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.755"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01784"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.616"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.75%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9009"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.00%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9984"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.28%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.986.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.27%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.96%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000121"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.82%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5096"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.35%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4065"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.36%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.960"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.61%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05763"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.17%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.745"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.03%  "
